$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new trade row (row 6) with the same shape/types as the existing rows.
$ws.Range("A6").Value = 9844.7900000000009
$ws.Range("B6").Value = 9893.27
$ws.Range("C6").Value = 307.20999999999998
$ws.Range("D6").Value = 308.70999999999998
$ws.Range("E6").Value = $true
$ws.Range("F6").Value = 0.49
$ws.Range("G6").Value = 42609.503946759258
$ws.Range("H6").Value = $false

# Copy the date/time number format from the row above (G5) onto G6 so it
# reuses the existing cell style instead of minting a new one.
$ws.Range("G5").Copy()
$ws.Range("G6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
